$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.030337353010300092
$ws.Range("B1").Value = -0.03033735442754968

$ws.Range("A2").Value = -0.055013197876232586
$ws.Range("B2").Value = 0.0550131964374138

$ws.Range("A3").Value = 0.069773215842998448
$ws.Range("B3").Value = -0.069773217209134986

$ws.Range("A4").Value = 0.062810646026219955
$ws.Range("B4").Value = -0.062810647434543274

$ws.Range("A5").Value = -0.036279071446375669
$ws.Range("B5").Value = 0.036279069987369426
